# Auto-generated Excel COM-interop script applying the diff to before.xlsx
# Sheets: 1=展览(Exhibition) 2=演出(Performance) 3=本地生活(Local Life) 4=全部类型(All Types)
$wb = $excel.ActiveWorkbook

# ---- Sheet1 (展览): update F-column 'want to go' counts ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 412
$ws1.Range("F4").Value = 145
$ws1.Range("F6").Value = 3702
$ws1.Range("F8").Value = 2497
$ws1.Range("F9").Value = 59
$ws1.Range("F10").Value = 2959
$ws1.Range("F12").Value = 525
$ws1.Range("F13").Value = 2256
$ws1.Range("F15").Value = 108
$ws1.Range("F16").Value = 38
$ws1.Range("F17").Value = 421
$ws1.Range("F19").Value = 182
$ws1.Range("F22").Value = 307
$ws1.Range("F24").Value = 1372
$ws1.Range("F27").Value = 117
$ws1.Range("F30").Value = 4084
$ws1.Range("F31").Value = 3648
$ws1.Range("F32").Value = 54
$ws1.Range("F34").Value = 1087
$ws1.Range("F35").Value = 440
$ws1.Range("F38").Value = 137
$ws1.Range("F39").Value = 118
$ws1.Range("F41").Value = 32
$ws1.Range("F43").Value = 41

# ---- Sheet2 (演出): update F-column 'want to go' counts ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 16
$ws2.Range("F16").Value = 192

# ---- Sheet4 (全部类型): update F-column 'want to go' counts ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 412
$ws4.Range("F7").Value = 145
$ws4.Range("F11").Value = 3702
$ws4.Range("F13").Value = 2497
$ws4.Range("F14").Value = 59
$ws4.Range("F15").Value = 2959
$ws4.Range("F16").Value = 525
$ws4.Range("F17").Value = 2256
$ws4.Range("F19").Value = 108
$ws4.Range("F20").Value = 38
$ws4.Range("F21").Value = 421
$ws4.Range("F23").Value = 182
$ws4.Range("F25").Value = 307
$ws4.Range("F27").Value = 1372
$ws4.Range("F33").Value = 4084
$ws4.Range("F34").Value = 3648
$ws4.Range("F35").Value = 54
$ws4.Range("F36").Value = 1087
$ws4.Range("F38").Value = 440
$ws4.Range("F44").Value = 137
$ws4.Range("F46").Value = 32
$ws4.Range("F48").Value = 41
$ws4.Range("F49").Value = 192

# ---- Sheet3 (本地生活): shift rows 3,4,5 content up into rows 2,3,4; update F values; delete old row 5 ----
$ws3 = $wb.Worksheets.Item(3)

# New row 2 (was row 3 data)
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = '2024-07-14'
$ws3.Range("B2").Style = "Normal"
$ws3.Range("C2").Value = '北京·排球少年!!垃圾场决战 主题咖啡厅'
$ws3.Range("D2").Value = '学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里'
$ws3.Range("E2").Value = '2024.07.14 00:00-08.25 23:59'
$ws3.Range("F2").Value = 1011
$ws3.Range("G2").Value = 10
$ws3.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=88981'
$ws3.Range("I2").Value = '//i1.hdslb.com/bfs/openplatform/202407/2T5mBYoB1720578883578.jpeg'

# New row 3 (was row 4 data)
$ws3.Range("B3").NumberFormat = "@"
$ws3.Range("B3").Value = '2024-07-17'
$ws3.Range("B3").Style = "Normal"
$ws3.Range("C3").Value = '北京·“狐妖小红娘”限时快闪店'
$ws3.Range("D3").Value = '王府井大街88号 北京王府井银泰in88购物中心'
$ws3.Range("E3").Value = '2024.07.17 10:00-10.31 22:00'
$ws3.Range("F3").Value = 137
$ws3.Range("G3").Value = 98
$ws3.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=89613'
$ws3.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202407/n3TXriJX1721203778030.jpeg'

# New row 4 (was row 5 data, but F updated 2186 -> 2193, G stays text)
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("B4").Value = '2024-08-17'
$ws3.Range("B4").Style = "Normal"
$ws3.Range("C4").Value = '北京·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅'
$ws3.Range("D4").Value = '酒仙桥恒通国际创新园C7栋012室 奇谷米·甜谷店（北京798店）'
$ws3.Range("E4").Value = '2024.08.17 00:00-10.09 23:59'
$ws3.Range("F4").Value = 2193
$ws3.Range("G4").Value = '已售罄'
$ws3.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=90435'
$ws3.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202408/0O7NsnOA1723429247959.png'

# Delete old row 5 entirely (shrinks dimension to A1:I4)
$ws3.Rows.Item(5).Delete()
